$d = $word.ActiveDocument

# The rPr child-element order inside several custom character styles in
# styles.xml violates the wml.xsd CT_RPr sequence (color must come after
# b/i, not before). Re-asserting Bold/Italic on the affected styles makes
# the writer re-serialize rPr in schema order without changing any actual
# formatting value.

$boldOnly = @("KeywordTok", "ImportTok", "ControlFlowTok", "AlertTok", "ErrorTok")
foreach ($styleName in $boldOnly) {
    $s = $d.Styles($styleName)
    $s.Font.Bold = 1
}

$italicOnly = @("CommentTok", "DocumentationTok")
foreach ($styleName in $italicOnly) {
    $s = $d.Styles($styleName)
    $s.Font.Italic = 1
}

$boldAndItalic = @("AnnotationTok", "CommentVarTok", "InformationTok", "WarningTok")
foreach ($styleName in $boldAndItalic) {
    $s = $d.Styles($styleName)
    $s.Font.Bold = 1
    $s.Font.Italic = 1
}
